$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended at row 19 (mirrors the "Append row" structure of row 18)
$ws.Range("A19").Value = "2025-05-01T11:53:37.511Z"
$ws.Range("B19").Value = "UNICEF"
$ws.Range("C19").Value = "C3"
$ws.Range("D19").Value = "الرحلة 2"
$ws.Range("E19").Value = "ايتا"
$ws.Range("F19").Value = "احمد"

# "1212" must stay text (not get auto-converted to a number) like the
# other numeric-looking strings in this sheet.
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "1212"
$ws.Range("G19").ClearFormats()

# H19 mirrors the empty string cells used throughout column H.
$ws.Range("H19").Font.Bold = $false
